$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Person rename: Jorrit de Boer -> Jelle Nauta (row 6)
# -----------------------------------------------------------------
$ws.Range("B6").Value = "Jelle"
$ws.Range("C6").Value = "Nauta"

# -----------------------------------------------------------------
# 2. Account row 14 (Acc_ad): mask password with formula, guard E14
# -----------------------------------------------------------------
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Formula = '=IF($B14="","","*****")'
$ws.Range("E14").Formula = '=IF($B14="","",$A$10)'

# -----------------------------------------------------------------
# 3. Account row 15 (Acc_rieks): mask password, guard E15
# -----------------------------------------------------------------
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Formula = '=IF($B15="","","*****")'
$ws.Range("E15").Formula = '=IF($B15="","",$A$10)'

# -----------------------------------------------------------------
# 4. Account row 16 (Acc_michiel): mask password, guard E16
# -----------------------------------------------------------------
$ws.Range("D16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Formula = '=IF($B16="","","*****")'
$ws.Range("E16").Formula = '=IF($B16="","",$A$10)'

# -----------------------------------------------------------------
# 5. Account row 17 (was Acc_jorrit -> Acc_jelle): userid + mask password, guard E17
# -----------------------------------------------------------------
$ws.Range("B17").Value = "jelle"
$ws.Range("D17").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Formula = '=IF($B17="","","*****")'
$ws.Range("E17").Formula = '=IF($B17="","",$A$10)'

# -----------------------------------------------------------------
# 6. Account row 18 (blank template row): add masked-password formula
#    to C18 and guard E18
# -----------------------------------------------------------------
$ws.Range("D18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Formula = '=IF($B18="","","*****")'
$ws.Range("E18").Formula = '=IF($B18="","",$A$10)'

# -----------------------------------------------------------------
# 7. Remove the blank spacer row (old row 19) so the [UIDs] block
#    and its data move up by one row.
# -----------------------------------------------------------------
$ws.Rows(19).Delete()

# -----------------------------------------------------------------
# 8. Fix up hyperlinks: the row delete above does not renumber the
#    hyperlink refs automatically, so rebuild them.
# -----------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A21"), "mailto:ms@g")
$ws.Hyperlinks.Add($ws.Range("F16"), "mailto:ms@google", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ms@google")
$ws.Hyperlinks.Add($ws.Range("A22"), "mailto:ms@li")
$ws.Hyperlinks.Add($ws.Range("A23"), "mailto:ms@tno")
$ws.Hyperlinks.Add($ws.Range("A24"), "mailto:rj@tno")
$ws.Hyperlinks.Add($ws.Range("A25"), "mailto:rj@g")

# -----------------------------------------------------------------
# 9. Selection cursor moves to A18 (matches the trimmed sheet)
# -----------------------------------------------------------------
$ws.Range("A18").Select()
